$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$Barr = New-Object 'object[,]' 24,1
$Barr[0,0] = 1.02
$Barr[1,0] = 1.02
$Barr[2,0] = 1.02
$Barr[3,0] = 1.02
$Barr[4,0] = 1.02
$Barr[5,0] = 1.02
$Barr[6,0] = 1.02
$Barr[7,0] = 1.02
$Barr[8,0] = 1.02
$Barr[9,0] = 1.02
$Barr[10,0] = 1.02
$Barr[11,0] = 1.02
$Barr[12,0] = 1.02
$Barr[13,0] = 1.02
$Barr[14,0] = 1.02
$Barr[15,0] = 1.02
$Barr[16,0] = 1.02
$Barr[17,0] = 1.02
$Barr[18,0] = 1.02
$Barr[19,0] = 1.02
$Barr[20,0] = 1.02
$Barr[21,0] = 1.02
$Barr[22,0] = 1.02
$Barr[23,0] = 1.02
$ws.Range("B2:B25").Value = $Barr

$Carr = New-Object 'object[,]' 24,1
$Carr[0,0] = 1.041250443480502
$Carr[1,0] = 1.042100009881367
$Carr[2,0] = 1.042650266558456
$Carr[3,0] = 1.042881719941534
$Carr[4,0] = 1.042920589297314
$Carr[5,0] = 1.04265335875873
$Carr[6,0] = 1.041537447693497
$Carr[7,0] = 1.039575198708455
$Carr[8,0] = 1.038269898168584
$Carr[9,0] = 1.037705387000796
$Carr[10,0] = 1.037495807962201
$Carr[11,0] = 1.037540758577333
$Carr[12,0] = 1.037688060970379
$Carr[13,0] = 1.037778832866623
$Carr[14,0] = 1.038307377634712
$Carr[15,0] = 1.038639106278784
$Carr[16,0] = 1.038832664708104
$Carr[17,0] = 1.038898674402667
$Carr[18,0] = 1.038603508022413
$Carr[19,0] = 1.037644681161511
$Carr[20,0] = 1.037042440044725
$Carr[21,0] = 1.037361641003946
$Carr[22,0] = 1.038619593136161
$Carr[23,0] = 1.040081988612102
$ws.Range("C2:C25").Value = $Carr

$Darr = New-Object 'object[,]' 24,1
$Darr[0,0] = 1.049925395660076
$Darr[1,0] = 1.050615647579132
$Darr[2,0] = 1.051062783848431
$Darr[3,0] = 1.051250877267311
$Darr[4,0] = 1.051282465811483
$Darr[5,0] = 1.051065296702979
$Darr[6,0] = 1.050158565585105
$Darr[7,0] = 1.048564671174197
$Darr[8,0] = 1.047504784894734
$Darr[9,0] = 1.047046506339157
$Darr[10,0] = 1.046876382020684
$Darr[11,0] = 1.046912869678346
$Darr[12,0] = 1.047032441744959
$Darr[13,0] = 1.047106127441845
$Darr[14,0] = 1.047535213367245
$Darr[15,0] = 1.047804545321622
$Darr[16,0] = 1.047961705582581
$Darr[17,0] = 1.048015303887091
$Darr[18,0] = 1.047775641972611
$Darr[19,0] = 1.046997227956937
$Darr[20,0] = 1.046508391875332
$Darr[21,0] = 1.046767477233619
$Darr[22,0] = 1.047788701955943
$Darr[23,0] = 1.048976260762944
$ws.Range("D2:D25").Value = $Darr

$Earr = New-Object 'object[,]' 24,1
$Earr[0,0] = 1.055781791390869
$Earr[1,0] = 1.056690544552168
$Earr[2,0] = 1.057279939554995
$Earr[3,0] = 1.057528047661681
$Earr[4,0] = 1.057569725196481
$Earr[5,0] = 1.057283253507147
$Earr[6,0] = 1.05608862385617
$Earr[7,0] = 1.053994111399813
$Earr[8,0] = 1.052604986483361
$Earr[9,0] = 1.052005213892084
$Earr[10,0] = 1.051782692998994
$Earr[11,0] = 1.051830412630427
$Earr[12,0] = 1.05198681490378
$Earr[13,0] = 1.052083214191601
$Earr[14,0] = 1.052644827967601
$Earr[15,0] = 1.052997577259532
$Earr[16,0] = 1.053203496527863
$Earr[17,0] = 1.05327373787641
$Earr[18,0] = 1.052959713356407
$Earr[19,0] = 1.051940751103515
$Earr[20,0] = 1.051301602083414
$Earr[21,0] = 1.051640283059723
$Earr[22,0] = 1.052976821909931
$Earr[23,0] = 1.054534328768894
$ws.Range("E2:E25").Value = $Earr

$Farr = New-Object 'object[,]' 24,1
$Farr[0,0] = 1.062229600228201
$Farr[1,0] = 1.063116080339917
$Farr[2,0] = 1.063690801365103
$Farr[3,0] = 1.063932677814979
$Farr[4,0] = 1.063973305348419
$Farr[5,0] = 1.063694032296029
$Farr[6,0] = 1.062528959603739
$Farr[7,0] = 1.060484526550035
$Farr[8,0] = 1.059127442748329
$Farr[9,0] = 1.058541225682678
$Farr[10,0] = 1.058323692096461
$Farr[11,0] = 1.058370344097062
$Farr[12,0] = 1.058523239922032
$Farr[13,0] = 1.058617472435713
$Farr[14,0] = 1.059166377889392
$Farr[15,0] = 1.059511070339617
$Farr[16,0] = 1.05971225945695
$Farr[17,0] = 1.059780882732713
$Farr[18,0] = 1.059474074037845
$Farr[19,0] = 1.058478210003213
$Farr[20,0] = 1.057853307099091
$Farr[21,0] = 1.058184462152987
$Farr[22,0] = 1.059490790654622
$Farr[23,0] = 1.061012034434144
$ws.Range("F2:F25").Value = $Farr

$Iarr = New-Object 'object[,]' 24,1
$Iarr[0,0] = 1.041695247779075
$Iarr[1,0] = 1.041893541597163
$Iarr[2,0] = 1.04202093751909
$Iarr[3,0] = 1.042074275420669
$Iarr[4,0] = 1.042083218207312
$Iarr[5,0] = 1.04202165108544
$Iarr[6,0] = 1.041762450836165
$Iarr[7,0] = 1.041298740511436
$Iarr[8,0] = 1.0409849561938
$Iarr[9,0] = 1.040847990939753
$Iarr[10,0] = 1.040796952153137
$Iarr[11,0] = 1.040807907547563
$Iarr[12,0] = 1.040843775395261
$Iarr[13,0] = 1.040865853075487
$Iarr[14,0] = 1.040994023121024
$Iarr[15,0] = 1.041074128217145
$Iarr[16,0] = 1.041120746509918
$Iarr[17,0] = 1.041136624202313
$Iarr[18,0] = 1.041065544622061
$Iarr[19,0] = 1.040833217724314
$Iarr[20,0] = 1.040686197253236
$Iarr[21,0] = 1.040764225163403
$Iarr[22,0] = 1.041069423505734
$Iarr[23,0] = 1.041419442214684
$ws.Range("I2:I25").Value = $Iarr

$Jarr = New-Object 'object[,]' 24,1
$Jarr[0,0] = 1.046332664107897
$Jarr[1,0] = 1.046828679993554
$Jarr[2,0] = 1.047149496247697
$Jarr[3,0] = 1.047284332902407
$Jarr[4,0] = 1.047306970512875
$Jarr[5,0] = 1.047151298078703
$Jarr[6,0] = 1.046500323112039
$Jarr[7,0] = 1.045352208022104
$Jarr[8,0] = 1.044586186186273
$Jarr[9,0] = 1.044254359883413
$Jarr[10,0] = 1.044131085893431
$Jarr[11,0] = 1.044157529420343
$Jarr[12,0] = 1.044244170396244
$Jarr[13,0] = 1.044297550288212
$Jarr[14,0] = 1.044608205686496
$Jarr[15,0] = 1.044803036618851
$Jarr[16,0] = 1.044916665217507
$Jarr[17,0] = 1.04495540739454
$Jarr[18,0] = 1.044782134438595
$Jarr[19,0] = 1.044218657301337
$Jarr[20,0] = 1.043864267940023
$Jarr[21,0] = 1.044052146377668
$Jarr[22,0] = 1.044791579274398
$Jarr[23,0] = 1.045649136252511
$ws.Range("J2:J25").Value = $Jarr

$Karr = New-Object 'object[,]' 24,1
$Karr[0,0] = 1.052681098155079
$Karr[1,0] = 1.053184008889538
$Karr[2,0] = 1.053509243766963
$Karr[3,0] = 1.053645927684643
$Karr[4,0] = 1.053668874850866
$Karr[5,0] = 1.053511070321862
$Karr[6,0] = 1.052851096121414
$Karr[7,0] = 1.051686798882094
$Karr[8,0] = 1.050909775620045
$Karr[9,0] = 1.050573135996622
$Karr[10,0] = 1.050448066751887
$Karr[11,0] = 1.050474895699354
$Karr[12,0] = 1.050562798263065
$Karr[13,0] = 1.050616954486388
$Karr[14,0] = 1.050932113525952
$Karr[15,0] = 1.05112975624481
$Karr[16,0] = 1.051245020116767
$Karr[17,0] = 1.051284319064026
$Karr[18,0] = 1.051108552881759
$Karr[19,0] = 1.05053691386486
$Karr[20,0] = 1.050177350188025
$Karr[21,0] = 1.050367975640147
$Karr[22,0] = 1.051118133823928
$Karr[23,0] = 1.051987948533672
$ws.Range("K2:K25").Value = $Karr

$Larr = New-Object 'object[,]' 24,1
$Larr[0,0] = 1.058521304632443
$Larr[1,0] = 1.059243324288861
$Larr[2,0] = 1.059711254413561
$Larr[3,0] = 1.059908146739156
$Larr[4,0] = 1.059941216016555
$Larr[5,0] = 1.059713884614647
$Larr[6,0] = 1.058765161446686
$Larr[7,0] = 1.057099086150051
$Larr[8,0] = 1.055992287907366
$Larr[9,0] = 1.055513980218585
$Larr[10,0] = 1.055336458668757
$Larr[11,0] = 1.055374531138439
$Larr[12,0] = 1.055499303301468
$Larr[13,0] = 1.05557619857265
$Larr[14,0] = 1.056024051615459
$Larr[15,0] = 1.056305231597767
$Larr[16,0] = 1.05646933003328
$Larr[17,0] = 1.05652529873297
$Larr[18,0] = 1.05627505423269
$Larr[19,0] = 1.055462557043551
$Larr[20,0] = 1.054952537202118
$Larr[21,0] = 1.055222829248992
$Larr[22,0] = 1.056288689804289
$Larr[23,0] = 1.057529122379199
$ws.Range("L2:L25").Value = $Larr

$Marr = New-Object 'object[,]' 24,1
$Marr[0,0] = 1.064951512862116
$Marr[1,0] = 1.065652587972155
$Marr[2,0] = 1.066106698466616
$Marr[3,0] = 1.066297717044608
$Marr[4,0] = 1.06632979634157
$Marr[5,0] = 1.066109250433813
$Marr[6,0] = 1.065188346658089
$Marr[7,0] = 1.063569248218407
$Marr[8,0] = 1.06249239498072
$Marr[9,0] = 1.062026728522499
$Marr[10,0] = 1.061853853661884
$Marr[11,0] = 1.061890931610081
$Marr[12,0] = 1.062012436698073
$Marr[13,0] = 1.062087312547209
$Marr[14,0] = 1.062523312883404
$Marr[15,0] = 1.062796970940003
$Marr[16,0] = 1.062956650589233
$Marr[17,0] = 1.063011107274637
$Marr[18,0] = 1.062767603851381
$Marr[19,0] = 1.061976653859986
$Marr[20,0] = 1.061479899112809
$Marr[21,0] = 1.061743185850449
$Marr[22,0] = 1.062780873390871
$Marr[23,0] = 1.063987381909924
$ws.Range("M2:M25").Value = $Marr

$Narr = New-Object 'object[,]' 24,1
$Narr[0,0] = 1.047818576506964
$Narr[1,0] = 1.048315296792073
$Narr[2,0] = 1.048636568642095
$Narr[3,0] = 1.048771596780321
$Narr[4,0] = 1.048794266538792
$Narr[5,0] = 1.048638373031907
$Narr[6,0] = 1.047986473606124
$Narr[7,0] = 1.046836728061069
$Narr[8,0] = 1.046069618386367
$Narr[9,0] = 1.045737320852097
$Narr[10,0] = 1.045613871798909
$Narr[11,0] = 1.045640352878663
$Narr[12,0] = 1.045727116894688
$Narr[13,0] = 1.045780572592224
$Narr[14,0] = 1.046091669156807
$Narr[15,0] = 1.046286776771433
$Narr[16,0] = 1.046400566735732
$Narr[17,0] = 1.046439363931101
$Narr[18,0] = 1.046265844907683
$Narr[19,0] = 1.045701567568259
$Narr[20,0] = 1.045346674933403
$Narr[21,0] = 1.045534820179979
$Narr[22,0] = 1.046275303156238
$Narr[23,0] = 1.047134077963622
$ws.Range("N2:N25").Value = $Narr
